{"js": "const body = context.document.body;\n\n// 1) \"Networked and Social Systems Engineering\" -> \"Networked & Social Systems Engineering\"\nlet results = body.search(\"and Social Systems Engineering\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"& Social Systems Engineering\", \"Replace\");\n}\n\n// 2) \"World\\u2019s largest network of local groups\" -> \"Worldwide network of in-person groups,\"\nresults = body.search(\"World\\u2019s largest network of local groups\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Worldwide network of in-person groups,\", \"Replace\");\n}\n\n// 3) \" 27 million members and 250,000 groups\" -> \" 28 million members and 260,000 groups\"\nresults = body.search(\" 27 million members and 250,000 groups\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\" 28 million members and 260,000 groups\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Networked and Social Systems Engineering\" -> \"Networked & Social Systems Engineering\"\n$r1 = $d.Content\n$r1.Find.Execute(\"and Social Systems Engineering\", $true, $false, $false, $false, $false, $true, 1, $false, \"& Social Systems Engineering\", 2)\n\n# 2) \"World's largest network of local groups\" -> \"Worldwide network of in-person groups,\"\n$r2 = $d.Content\n$r2.Find.Execute(\"World\u2019s largest network of local groups\", $true, $false, $false, $false, $false, $true, 1, $false, \"Worldwide network of in-person groups,\", 2)\n\n# 3) \" 27 million members and 250,000 groups\" -> \" 28 million members and 260,000 groups\"\n$r3 = $d.Content\n$r3.Find.Execute(\" 27 million members and 250,000 groups\", $true, $false, $false, $false, $false, $true, 1, $false, \" 28 million members and 260,000 groups\", 2)\n"}
